# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1. Update the "last updated" timestamp label in A1
$ws.Range("A1").Value = "Datos actualizados a 15 de Agosto de 2020 a las 08:04"

# 2. Update country data rows (India, Pakistan, Afganistan, Uzbekistan)
# Row 6: India
$ws.Range("B6").Value = 2527308
$ws.Range("C6").Value = 2086
$ws.Range("D6").Value = 1809542
$ws.Range("E6").Value = 668618
$ws.Range("G6").Value = 14
$ws.Range("H6").Value = 49148

# Row 17: Pakistan
$ws.Range("B17").Value = 288047
$ws.Range("C17").Value = 747
$ws.Range("D17").Value = 265624
$ws.Range("E17").Value = 16261
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = 6162

# Row 60: Afganistan
$ws.Range("B60").Value = 37550
$ws.Range("C60").Value = 119
$ws.Range("D60").Value = 27166
$ws.Range("E60").Value = 9015
$ws.Range("G60").Value = 6
$ws.Range("H60").Value = 1369

# Row 62: Uzbekistan
$ws.Range("B62").Value = 34017
$ws.Range("C62").Value = 196
$ws.Range("D62").Value = 28661
$ws.Range("E62").Value = 5135
$ws.Range("G62").Value = 1
$ws.Range("H62").Value = 221

# 3. Rows 213/214: Montserrat and Islas Malvinas swap ranking order
#    Row 213 now holds Montserrat's data, row 214 now holds Islas Malvinas's data
$ws.Range("A213").Value = "Montserrat"
$ws.Range("D213").Value = 12
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0
